# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the "7a0011b4-1f4f-415f-8bde-5f3a69dbe1b5" item.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth setter stores a value ~0.8333 chars wider than the
# number assigned (character-padding quirk), so pre-compensate to land on
# an on-disk column width of exactly 40.
$targetColWidth = 40 - 0.8333333333333334

# --- Overview sheet: status text for the 7a0011b4 row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: Status (column C) + Error Detail (column P) for the 7a0011b4 row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: gxyl3d4j.zke is different with handoff file name: 7a0011b4-1f4f-415f-8bde-5f3a69dbe1b5.ff5cb0c85822e0e75887d3fdf20608c722064911.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet: Status (column C) + Error Detail (column P) for the 7a0011b4 row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: gxyl3d4j.zke is different with handoff file name: 7a0011b4-1f4f-415f-8bde-5f3a69dbe1b5.ff5cb0c85822e0e75887d3fdf20608c722064911.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
